# Update cryptocurrency price/volume snapshot cells (columns D and E)
# on Sheet1. Column D ("Price") is stored as text in the workbook (it mixes
# thousands-dot-formatted values like "30.225.56" with plain decimals like
# "0.9998"), so purely-numeric-looking price strings are entered with a
# leading apostrophe to force Excel to keep them as text instead of
# auto-converting them to numbers (this also matches how the sheet was
# originally authored).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.225.56"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.859.75"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'236.71"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4675"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.2867"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").Value = "'0.06541"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +5.56%  "
$ws.Range("D11").Value = "'0.07927"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "'97.74"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "1.865.80"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "'5.181"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "'0.6805"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "'267.92"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "30.218.32"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'13.77"
$ws.Range("E18").Value = "  +8.81%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "2.111.73"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'5.329"
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'6.206"
$ws.Range("D25").Value = "'167.38"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").Value = "'9.222"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").Value = "'1.386"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("D30").Value = "'0.09895"
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").Value = "'4.390"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "'4.070"
$ws.Range("D34").Value = "'0.04706"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").Value = "'0.7050"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'0.01891"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "'2.633"
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("D40").Value = "'6.268"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "'74.37"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").Value = "'1.940"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'0.8481"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'0.4171"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'0.9991"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'103.54"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "'967.19"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("D48").Value = "'7.170"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").Value = "'9.243"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "'34.17"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'0.05653"
$ws.Range("E51").Value = "  +0.42%  "
